$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exceptions")

# Add a new row describing the new 0x0005 error first.
$ws.Range("A7").Value = "0x0005"
$ws.Range("B7").Value = "TopModel.cs"

# Update existing row 6, column C: shorten the message text,
# the error code itself now lives in its own row below.
$ws.Range("C6").Value = "Failed to ReadFromTable,"

$ws.Range("C7").Value = "Try Removing Top Failed"

# Keep the active selection in sync with the new last cell, like Excel would
# after the edit.
$ws.Range("C7").Select()
